$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.378.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.72"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6305"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.23%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07544"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2956"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.43"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.12%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07702"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.42%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.81"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.995"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6832"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.31%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001003"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.67%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.80%  "

# Row 17
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.117.88"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.137"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.406.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.557"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.12%  "

# Row 24
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.950"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.000"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.64%  "

# Row 27
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1396"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.49%  "

# Row 28
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.373"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.95%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.65"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.29%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.464"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.72%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.264"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.90%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05680"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.31%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.121"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.30%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.020"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.847"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.36%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.154"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.29%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7124"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.02%  "

# Row 38
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.591"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.12%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.261.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.87%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01814"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.02%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.778"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.70%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9080"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.12%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.167"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.38%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.076"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.45%  "

# Row 48
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4043"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.00%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.103"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.684"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.80%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1120"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.07%  "
